$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append the new log entry as row 55
$ws.Range("A55").Value = "Sollicitatie marketingfunctie"
$ws.Range("B55").Value = "mailmind.test@zohomail.eu"
$ws.Range("C55").Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$ws.Range("D55").Value = "Overig"
$ws.Range("F55").Value = "2025-06-17 22:30:06"
$ws.Range("G55").Value = "Nee"

# Extend the conditional formatting ranges to include the new row
$fcsD = $ws.Range("D2:D54").FormatConditions
$fcsD.Item(1).ModifyAppliesToRange($ws.Range("D2:D55"))

$fcsG = $ws.Range("G2:G54").FormatConditions
$fcsG.Item(1).ModifyAppliesToRange($ws.Range("G2:G55"))

# Update the Dashboard summary count for the "Overig" category
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 15
